$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 144541.92
$ws.Range("J17").Value = 144541.92
$ws.Range("L17").Value = 433625.76
$ws.Range("N17").Value = -433961.76

$ws.Range("H19").Value = 1318.0625
$ws.Range("J19").Value = 864.4286
$ws.Range("L19").Value = 864.4286
$ws.Range("N19").Value = -1214.4286

$ws.Range("H132").Value = 1151.4375
$ws.Range("I132").Value = 1197.8966
$ws.Range("K132").Value = 3593.6898
$ws.Range("M132").Value = -1063.6898

$ws.Range("H135").Value = 88235816
$ws.Range("I135").Value = 31250556
$ws.Range("J135").Value = 1000000000
$ws.Range("K135").Value = 281255004
$ws.Range("L135").Value = 9000000000
$ws.Range("M135").Value = -281252469
$ws.Range("N135").Value = -9000005070

$ws.Range("H137").Value = 3385.6155
$ws.Range("I137").Value = 3084.4167
$ws.Range("K137").Value = 9253.250100000001
$ws.Range("M137").Value = -6703.250100000001

$ws.Range("H141").Value = 3837.25
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4052.1365
$ws.Range("I32").Value = 2982.2974
$ws.Range("J32").Value = 9707
$ws.Range("K32").Value = 2982.2974
$ws.Range("L32").Value = 9707
$ws.Range("M32").Value = -2695.2974
$ws.Range("N32").Value = -10281

$ws.Range("H45").Value = 1974.5
$ws.Range("I45").Value = 1650
$ws.Range("K45").Value = 1650
$ws.Range("M45").Value = -1273

$ws.Range("H61").Value = 200010620
$ws.Range("I61").Value = 200010620
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 200010620
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -200010408
$ws.Range("N61").ClearContents()

$ws.Range("H122").Value = 3899.4546
$ws.Range("I122").Value = 2321.6667
$ws.Range("K122").Value = 6965.000100000001
$ws.Range("M122").Value = -4515.000100000001

$ws.Range("H136").Value = 200010620
$ws.Range("I136").Value = 200010620
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 600031860
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -600029310
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H102").Value = 3999.6667
$ws.Range("I102").Value = 3999.6667
$ws.Range("K102").Value = 3999.6667
$ws.Range("M102").Value = -754.6667000000002

$ws.Range("H134").Value = 100002856
$ws.Range("I134").Value = 125002350
$ws.Range("J134").Value = 4900
$ws.Range("K134").Value = 375007050
$ws.Range("L134").Value = 14700
$ws.Range("M134").Value = -375004515
$ws.Range("N134").Value = -19770

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1557480.8
$ws.Range("I16").Value = 2718841.2
$ws.Range("K16").Value = 2718841.2
$ws.Range("M16").Value = -2718554.2

$ws.Range("H31").Value = 11857.038
$ws.Range("I31").Value = 6080.75
$ws.Range("K31").Value = 6080.75
$ws.Range("M31").Value = -5785.75

$ws.Range("H34").Value = 11857.038
$ws.Range("I34").Value = 6080.75
$ws.Range("K34").Value = 6080.75
$ws.Range("M34").Value = -5878.75

$ws.Range("H86").Value = 3981.2
$ws.Range("I86").Value = 3803.6667
$ws.Range("J86").Value = 4247.5
$ws.Range("K86").Value = 3803.6667
$ws.Range("L86").Value = 4247.5
$ws.Range("M86").Value = -2680.6667
$ws.Range("N86").Value = -6493.5

$ws.Range("H89").Value = 3981.2
$ws.Range("I89").Value = 3803.6667
$ws.Range("J89").Value = 4247.5
$ws.Range("K89").Value = 19018.3335
$ws.Range("L89").Value = 21237.5
$ws.Range("M89").Value = -13402.3335
$ws.Range("N89").Value = -32469.5

$ws.Range("H113").Value = 1557480.8
$ws.Range("I113").Value = 2718841.2
$ws.Range("K113").Value = 2718841.2
$ws.Range("M113").Value = -2716671.2

$ws.Range("H141").Value = 352861
$ws.Range("J141").Value = 388984
$ws.Range("L141").Value = 388984
$ws.Range("N141").Value = -399344

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 2122.25
$ws.Range("I63").Value = 2122.25
$ws.Range("K63").Value = 6366.75
$ws.Range("M63").Value = -5617.75

$ws.Range("H66").Value = 2122.25
$ws.Range("I66").Value = 2122.25
$ws.Range("K66").Value = 19100.25
$ws.Range("M66").Value = -15356.25

$ws.Range("H131").Value = 2029.1
$ws.Range("I131").Value = 2300
$ws.Range("K131").Value = 6900
$ws.Range("M131").Value = -1860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 48999.332
$ws.Range("I18").Value = 48999
$ws.Range("K18").Value = 48999
$ws.Range("M18").Value = -48706

$ws.Range("H122").Value = 4221.727
$ws.Range("I122").Value = 1910
$ws.Range("J122").Value = 5542.7144
$ws.Range("K122").Value = 5730
$ws.Range("L122").Value = 16628.1432
$ws.Range("M122").Value = -3280
$ws.Range("N122").Value = -21528.1432

$ws.Range("H126").Value = 6360.4287
$ws.Range("I126").Value = 6360.4287
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 19081.2861
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -16611.2861
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2403.3
$ws.Range("I40").Value = 2348.111
$ws.Range("K40").Value = 2348.111
$ws.Range("M40").Value = -2212.111

$ws.Range("H100").Value = 13309529
$ws.Range("I100").Value = 22178436
$ws.Range("K100").Value = 22178436
$ws.Range("M100").Value = -22177895

$ws.Range("H122").Value = 8595.223
$ws.Range("I122").Value = 8595.223
$ws.Range("K122").Value = 25785.669
$ws.Range("M122").Value = -23335.669

$ws.Range("H132").Value = 31277600
$ws.Range("I132").Value = 35745500
$ws.Range("K132").Value = 107236500
$ws.Range("M132").Value = -107233970

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 9950
$ws.Range("J7").Value = 9950
$ws.Range("L7").Value = 9950
$ws.Range("N7").Value = -10176

$ws.Range("H107").Value = 1511.2
$ws.Range("I107").Value = 1246.7142
$ws.Range("K107").Value = 3740.1426
$ws.Range("M107").Value = -1820.1426

$ws.Range("H122").Value = 5419.7144
$ws.Range("I122").Value = 5883.1665
$ws.Range("K122").Value = 17649.4995
$ws.Range("M122").Value = -15199.4995

$ws.Range("H132").Value = 12825496
$ws.Range("I132").Value = 20002398
$ws.Range("J132").Value = 9598.786
$ws.Range("K132").Value = 60007194
$ws.Range("L132").Value = 28796.358
$ws.Range("M132").Value = -60004664
$ws.Range("N132").Value = -33856.358

$ws.Range("H136").Value = 21743510
$ws.Range("I136").Value = 25004260
$ws.Range("J136").Value = 5186.3335
$ws.Range("K136").Value = 75012780
$ws.Range("L136").Value = 15559.0005
$ws.Range("M136").Value = -75010230
$ws.Range("N136").Value = -20659.0005

